$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S2").Value = 3.35
$ws.Range("T2").Value = 1.72
$ws.Range("W2").Value = 1.45
$ws.Range("F3").Value = 1.84
$ws.Range("G3").Value = 1.86
$ws.Range("H3").Value = 4.5
$ws.Range("L3").Value = 1.27
$ws.Range("N3").Value = 5.7
$ws.Range("O3").Value = 1.19
$ws.Range("P3").Value = 2.58
$ws.Range("Q3").Value = 1.6
$ws.Range("R3").Value = 1.64
$ws.Range("S3").Value = 2.44
$ws.Range("T3").Value = 1.59
$ws.Range("U3").Value = 2.58
$ws.Range("V3").Value = 1.27
$ws.Range("W3").Value = 2.18
$ws.Range("X3").Value = 23
$ws.Range("Z3").Value = 38
$ws.Range("AA3").Value = 100
$ws.Range("AB3").Value = 13.5
$ws.Range("AC3").Value = 9.800000000000001
$ws.Range("AD3").Value = 17.5
$ws.Range("AE3").Value = 46
$ws.Range("AF3").Value = 14
$ws.Range("AG3").Value = 9.800000000000001
$ws.Range("AI3").Value = 46
$ws.Range("AK3").Value = 16
$ws.Range("AL3").Value = 26
$ws.Range("AM3").Value = 65
$ws.Range("AN3").Value = 8.199999999999999
$ws.Range("AO3").Value = 36
$ws.Range("G4").Value = 12.5
$ws.Range("H4").Value = 1.34
$ws.Range("J4").Value = 5.3
$ws.Range("L4").Value = 1.01
$ws.Range("Q4").Value = 1.6
$ws.Range("U4").Value = 1.81
$ws.Range("V4").Value = 3.5
$ws.Range("W4").Value = 1.08
$ws.Range("AC4").Value = 14
$ws.Range("AD4").Value = 11.5
$ws.Range("F5").Value = 1.81
$ws.Range("G5").Value = 1.93
$ws.Range("I5").Value = 4.7
$ws.Range("J5").Value = 3.85
$ws.Range("K5").Value = 4.6
$ws.Range("L5").Value = 1.01
$ws.Range("M5").Value = 1.04
$ws.Range("N5").Value = 1.01
$ws.Range("O5").Value = 1.21
$ws.Range("R5").Value = 1.19
$ws.Range("S5").Value = 1.01
$ws.Range("T5").Value = 1.01
$ws.Range("U5").Value = 2.1
$ws.Range("V5").Value = 1.27
$ws.Range("W5").Value = 2.06
$ws.Range("X5").Value = 1000
$ws.Range("Y5").Value = 30
$ws.Range("Z5").Value = 1000
$ws.Range("AA5").Value = 1000
$ws.Range("AB5").Value = 17.5
$ws.Range("AC5").Value = 14.5
$ws.Range("AD5").Value = 26
$ws.Range("AE5").Value = 1000
$ws.Range("AF5").Value = 19
$ws.Range("AG5").Value = 15.5
$ws.Range("AH5").Value = 25
$ws.Range("AI5").Value = 1000
$ws.Range("AJ5").Value = 30
$ws.Range("AK5").Value = 26
$ws.Range("AL5").Value = 44
$ws.Range("AM5").Value = 1000
$ws.Range("AN5").Value = 1000
$ws.Range("AO5").Value = 1000
$ws.Range("J6").Value = 5.9
$ws.Range("AD6").Value = 60
$ws.Range("F7").Value = 3.7
$ws.Range("P7").Value = 2.38
$ws.Range("AC7").Value = 9.6
$ws.Range("AH7").Value = 19
$ws.Range("AO7").Value = 12.5
$ws.Range("H8").Value = 6
$ws.Range("J8").Value = 3.95
$ws.Range("N8").Value = 3.7
$ws.Range("P8").Value = 1.92
$ws.Range("T8").Value = 2.02
$ws.Range("AC8").Value = 8.800000000000001
$ws.Range("AH8").Value = 23
$ws.Range("F9").Value = 1.94
$ws.Range("G9").Value = 2.1
$ws.Range("H9").Value = 3.85
$ws.Range("I9").Value = 4.5
$ws.Range("J9").Value = 3.6
$ws.Range("K9").Value = 4.2
$ws.Range("P9").Value = 2
$ws.Range("Q9").Value = 1.82
